# PlayerPerformance_3877.xlsx edit script
# Adds a "Player Info" sheet (new first sheet) and an "ODI Batting Extra" sheet
# (new last sheet), renames MATCH_CARD_LINK -> MATCH_CODE on the existing
# "ODI Batting"/"ODI Bowling" sheets and rewrites those URL columns down to
# the bare match code.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell forcing it to be stored as TEXT even when
# it looks numeric (so "3310" stays a string instead of becoming 3310).
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------------
# 1) Add the two brand new sheets at the END of the workbook first, so that
#    the existing "ODI Batting" / "ODI Bowling" sheets keep their original
#    positions (and so variables referencing them by name stay valid) while
#    we populate the new sheets. We reorder everything at the very end.
# ---------------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$playerInfoWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$playerInfoWs.Name = "Player Info"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtraWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$battingExtraWs.Name = "ODI Batting Extra"

# ---------------------------------------------------------------------------
# 2) Populate "Player Info"
# ---------------------------------------------------------------------------
$playerInfoWs = $wb.Worksheets.Item("Player Info")

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfoWs.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

Set-TextValue $playerInfoWs.Cells.Item(2, 1) "3877"
$playerInfoWs.Cells.Item(2, 2).Value = "Gulbadin Naib"
$playerInfoWs.Cells.Item(2, 3).Value = "Right Handed"
$playerInfoWs.Cells.Item(2, 4).Value = "Right Arm Medium Fast"

# ---------------------------------------------------------------------------
# 3) Fix up "ODI Batting" (MATCH_CARD_LINK -> MATCH_CODE, header + values)
# ---------------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"

$lastRowBatting = $wsBatting.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRowBatting; $r++) {
    $cell = $wsBatting.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $code = $val -replace '^.*MatchCode=', ''
        Set-TextValue $cell $code
    }
}

# ---------------------------------------------------------------------------
# 4) Fix up "ODI Bowling" (MATCH_CARD_LINK -> MATCH_CODE, header + values)
# ---------------------------------------------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsBowling.Cells.Item(1, 2).Value = "MATCH_CODE"

$lastRowBowling = $wsBowling.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRowBowling; $r++) {
    $cell = $wsBowling.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $code = $val -replace '^.*MatchCode=', ''
        Set-TextValue $cell $code
    }
}

# ---------------------------------------------------------------------------
# 5) Populate "ODI Batting Extra"
# ---------------------------------------------------------------------------
$battingExtraWs = $wb.Worksheets.Item("ODI Batting Extra")

$beHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $beHeaders.Length; $c++) {
    $cell = $battingExtraWs.Cells.Item(1, $c)
    $cell.Value = $beHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# MATCH_CODE, BATTING_POSITION (numeric or $null), NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$beRows = @(
    @("4306", 6, "4", "1", "14.98%", "NO"),
    @("4309", 6, "2", "0", "15.13%", "NO"),
    @("4315", 5, "1", "0", "2.33%", "NO"),
    @("4323", 8, "0", "0", "4.00%", "NO"),
    @("4326", $null, $null, $null, $null, "NO"),
    @("4332", $null, $null, $null, $null, "NO"),
    @("4335", 1, "3", "0", "23.50%", "NO"),
    @("4340", 2, "3", "0", "6.61%", "NO"),
    @("4348", 1, "1", "0", "1.74%", "NO"),
    @("4377", 8, "3", "0", "8.76%", "NO"),
    @("4444", 6, "0", "0", "0.70%", "NO"),
    @("4446", 8, $null, $null, $null, "NO"),
    @("4448", 7, "3", "1", "13.53%", "NO"),
    @("4525", $null, $null, $null, $null, "NO"),
    @("4528", 6, "0", "0", "0.84%", "NO"),
    @("4537", 7, "0", "1", "7.91%", "NO"),
    @("4539", 7, $null, $null, $null, "NO"),
    @("4671", $null, $null, $null, $null, "NO"),
    @("4674", 6, "0", "0", "2.63%", "NO"),
    @("4675", $null, $null, $null, $null, "NO")
)

$rowNum = 2
foreach ($row in $beRows) {
    Set-TextValue $battingExtraWs.Cells.Item($rowNum, 1) $row[0]
    if ($row[1] -ne $null) {
        $battingExtraWs.Cells.Item($rowNum, 2).Value = $row[1]
    }
    if ($row[2] -ne $null) {
        Set-TextValue $battingExtraWs.Cells.Item($rowNum, 3) $row[2]
    }
    if ($row[3] -ne $null) {
        Set-TextValue $battingExtraWs.Cells.Item($rowNum, 4) $row[3]
    }
    if ($row[4] -ne $null) {
        Set-TextValue $battingExtraWs.Cells.Item($rowNum, 5) $row[4]
    }
    Set-TextValue $battingExtraWs.Cells.Item($rowNum, 6) $row[5]
    $rowNum++
}

# ---------------------------------------------------------------------------
# 6) Re-order sheets: Player Info, ODI Batting, ODI Bowling, ODI Batting Extra
# ---------------------------------------------------------------------------
$playerInfoWs = $wb.Worksheets.Item("Player Info")
$firstSheet = $wb.Worksheets.Item(1)
$playerInfoWs.Move($firstSheet)

$wb.Worksheets.Item(1).Select()
